$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("d-dataSheet").Delete()
$wb.Worksheets.Item("d-dev").Delete()
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
